$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update URL, Version, Date, Publisher ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-id"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet: clear the stray duplicated constraint text on the
#     "Extension" summary row (row 2, column AI = Constraint(s)) ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

Write-Output "edits applied"
